# edit.ps1
# Applies the commit "Part 3 finish - need DEBUG" to before.docx:
#  1. Append a "." run to the "In File" paragraph, then insert a brand-new
#     paragraph right after it ("Evaluate the labels in a naive approach...").
#  2. Split the "We assume that..." paragraph's single run into three runs,
#     wrapping "is" with gramStart/gramEnd proofErr markers.
#  3. After the final (image) paragraph, insert a blank paragraph, the
#     "Part C: ..." heading, the "5. Implement the score method..." body
#     paragraph and a trailing blank paragraph.

$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Edit 1a: "In File" -> "In File."  (new run appended in the same paragraph,
# preserving the paragraph's own identity/attributes).
# ---------------------------------------------------------------------
$pInFile = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "In File") {
        $pInFile = $p
        break
    }
}

$rr1 = $d.Range($pInFile.Range.Start, $pInFile.Range.End - 1)
$xml1 = "<w:p $ns><w:pPr><w:bidi w:val=`"0`"/></w:pPr>" +
        "<w:r><w:t>In File</w:t></w:r>" +
        "<w:r><w:t>.</w:t></w:r>" +
        "</w:p>"
$rr1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Edit 1b: insert the brand-new "Evaluate the labels..." paragraph right
# after the "In File." paragraph.
# ---------------------------------------------------------------------
$rAfter = $pInFile.Range
$rAfter.Collapse(0)
$rAfter.InsertParagraphAfter()

$pNew = $pInFile.Next()
$rr2 = $d.Range($pNew.Range.Start, $pNew.Range.End)
$xml2 = "<w:p $ns><w:pPr><w:bidi w:val=`"0`"/></w:pPr>" +
        "<w:r><w:t>Evaluate the labels in a naive approach. Each value in the</w:t></w:r>" +
        "<w:r><w:br/><w:t xml:space=`"preserve`">result tensor should contain the disparity matching minimal </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>ssd</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> (sum of</w:t></w:r>" +
        "<w:r><w:br/><w:t>squared difference).</w:t></w:r>" +
        "</w:p>"
$rr2.InsertXML($xml2)

# ---------------------------------------------------------------------
# Edit 2: split "We assume that ... figure is good but ..." into three
# runs, marking "is" with gramStart/gramEnd proofErr tags.
# ---------------------------------------------------------------------
$pAssume = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("We assume that the disparity")) {
        $pAssume = $p
        break
    }
}

$naive = "na" + [char]0x00EF + "ve"
$rr3 = $d.Range($pAssume.Range.Start, $pAssume.Range.End - 1)
$xml3 = "<w:p $ns><w:pPr><w:bidi w:val=`"0`"/></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">We assume that the disparity should be the same for all the pixels. But, in the result we can see that when the colors are less likely to the colors around, the results for the figure </w:t></w:r>" +
        "<w:proofErr w:type=`"gramStart`"/>" +
        "<w:r><w:t>is</w:t></w:r>" +
        "<w:proofErr w:type=`"gramEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> good but wen we are having smooth colors change, the $naive depth map in noisy.</w:t></w:r>" +
        "</w:p>"
$rr3.InsertXML($xml3)

# ---------------------------------------------------------------------
# Edit 3: after the final (image) paragraph, append a blank paragraph,
# the "Part C: ..." heading, the "5. Implement..." body paragraph, and a
# trailing blank paragraph.
# ---------------------------------------------------------------------
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$rLastEnd = $pLast.Range
$rLastEnd.Collapse(0)
$rLastEnd.InsertParagraphAfter()

$pTail = $d.Paragraphs($d.Paragraphs.Count)
$rr4 = $d.Range($pTail.Range.Start, $pTail.Range.End)
$xml4 = "<w:p $ns><w:pPr><w:bidi w:val=`"0`"/></w:pPr></w:p>" +
        "<w:p $ns><w:pPr><w:pStyle w:val=`"1`"/><w:bidi w:val=`"0`"/></w:pPr>" +
        "<w:r><w:rPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"FFFFFF`"/></w:rPr>" +
        "<w:lastRenderedPageBreak/>" +
        "<w:t>Part C: Depth Map Smoothing using Dynamic Programming</w:t></w:r>" +
        "</w:p>" +
        "<w:p $ns><w:pPr><w:pStyle w:val=`"2`"/><w:bidi w:val=`"0`"/></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">5. </w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`">Implement the score method for a single slice of the </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/>" +
        "<w:r><w:t>ssdd</w:t></w:r>" +
        "<w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> tensor, using Dynamic Programming.</w:t></w:r>" +
        "</w:p>" +
        "<w:p $ns><w:pPr><w:bidi w:val=`"0`"/></w:pPr></w:p>"
$rr4.InsertXML($xml4)
